# Insert two new columns before column AE, pushing the existing
# AE:AO ("Note" .. "MILOFtime_taken") block to AG:AQ.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AE1:AF1").EntireColumn.Insert()

# New header cells for the two inserted "Unnamed: 0...." id columns,
# continuing the existing dot-count naming pattern (AD1 has 28 dots).
$ws.Range("AE1").Value = "Unnamed: 0.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1"
$ws.Range("AF1").Value = "Unnamed: 0.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1.1"

# The two new id columns mirror column AD's running row index (0..6)
# for each data row, same as the other "Unnamed: 0..." index columns.
for ($row = 2; $row -le 8; $row++) {
    $idValue = $ws.Cells.Item($row, 30).Value2
    $ws.Cells.Item($row, 31).Value = $idValue
    $ws.Cells.Item($row, 32).Value = $idValue
}

# Re-run MILOF for the two datasets whose results changed
# (ambient_temperature_system_failure.csv and
# cpu_utilization_asg_misconfiguration.csv): new identified-discord
# list, new best params, new time taken.
$ws.Range("AN2").Value = "[3213, 3637, 6012]"
$ws.Range("AP2").Value = "{'Numk': 17, 'KPar': 4, 'Bucket_index': 500}"
$ws.Range("AQ2").Value = 178.6288073339965

$ws.Range("AN3").Value = "[16727, 17627, 17951]"
$ws.Range("AP3").Value = "{'Numk': 23, 'KPar': 14, 'Bucket_index': 500}"
$ws.Range("AQ3").Value = 430.6199549960438
